# Update Seraph_Profits sheets with refreshed market-price derived values
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 2000
$ws.Range("J49").Value = 2000
$ws.Range("L49").Value = 6000
$ws.Range("N49").Value = -6272

$ws.Range("H81").Value = 38000
$ws.Range("J81").Value = 38000
$ws.Range("L81").Value = 38000
$ws.Range("N81").Value = -39996

$ws.Range("H84").Value = 38000
$ws.Range("J84").Value = 38000
$ws.Range("L84").Value = 114000
$ws.Range("N84").Value = -123984

$ws.Range("H96").Value = 3439.7778
$ws.Range("I96").Value = 3080
$ws.Range("J96").Value = 3727.6
$ws.Range("K96").Value = 9240
$ws.Range("L96").Value = 11182.8
$ws.Range("M96").Value = -7867
$ws.Range("N96").Value = -13928.8

$ws.Range("H100").Value = 1397.3572
$ws.Range("I100").Value = 1022.4
$ws.Range("J100").Value = 2334.75
$ws.Range("K100").Value = 1022.4
$ws.Range("L100").Value = 2334.75
$ws.Range("M100").Value = -481.4
$ws.Range("N100").Value = -3416.75

$ws.Range("H107").Value = 271.3846
$ws.Range("I107").Value = 271.3846
$ws.Range("K107").Value = 271.3846
$ws.Range("M107").Value = 1648.6154

$ws.Range("H113").Value = 1832.8334
$ws.Range("I113").Value = 1819.4
$ws.Range("K113").Value = 1819.4
$ws.Range("M113").Value = 1434.6

$ws.Range("H138").Value = 4185.864
$ws.Range("I138").Value = 1534.125
$ws.Range("J138").Value = 4775.1387
$ws.Range("K138").Value = 4602.375
$ws.Range("L138").Value = 14325.4161
$ws.Range("M138").Value = 537.625
$ws.Range("N138").Value = -24605.4161

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11178.333
$ws.Range("I32").Value = 6546.1763
$ws.Range("K32").Value = 6546.1763
$ws.Range("M32").Value = -6259.1763

$ws.Range("H61").Value = 4596.5
$ws.Range("J61").Value = 4596.5
$ws.Range("L61").Value = 4596.5
$ws.Range("N61").Value = -5020.5

$ws.Range("H63").Value = 5020.7334
$ws.Range("J63").Value = 6869.625
$ws.Range("L63").Value = 6869.625
$ws.Range("N63").Value = -8241.625

$ws.Range("H66").Value = 5020.7334
$ws.Range("J66").Value = 6869.625
$ws.Range("L66").Value = 34348.125
$ws.Range("N66").Value = -41212.125

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws.Range("H110").Value = 17397.2
$ws.Range("I110").Value = 19425.334
$ws.Range("K110").Value = 19425.334
$ws.Range("M110").Value = -17380.334

$ws.Range("H136").Value = 4596.5
$ws.Range("J136").Value = 4596.5
$ws.Range("L136").Value = 13789.5
$ws.Range("N136").Value = -18889.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 761.5357
$ws.Range("I94").Value = 801.0769
$ws.Range("K94").Value = 801.0769
$ws.Range("M94").Value = -350.0769

$ws.Range("H105").Value = 5411.6875
$ws.Range("I105").Value = 3772.5334
$ws.Range("K105").Value = 3772.5334
$ws.Range("M105").Value = -2025.5334

$ws.Range("H107").Value = 659.3043
$ws.Range("I107").Value = 658.3333
$ws.Range("K107").Value = 658.3333
$ws.Range("M107").Value = 1261.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 29711.5
$ws.Range("J51").Value = 49999
$ws.Range("L51").Value = 49999
$ws.Range("N51").Value = -51471

$ws.Range("H52").Value = 86333
$ws.Range("J52").Value = 86333
$ws.Range("L52").Value = 86333
$ws.Range("N52").Value = -86921

$ws.Range("H61").Value = 29711.5
$ws.Range("J61").Value = 49999
$ws.Range("L61").Value = 49999
$ws.Range("N61").Value = -50695

$ws.Range("H105").Value = 1233.4
$ws.Range("I105").Value = 1233.4
$ws.Range("K105").Value = 1233.4
$ws.Range("M105").Value = 513.5999999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 33333416
$ws.Range("J7").Value = 45
$ws.Range("L7").Value = 135
$ws.Range("N7").Value = -359

$ws.Range("H34").Value = 2002.4166
$ws.Range("I34").Value = 907.75
$ws.Range("J34").Value = 2549.75
$ws.Range("K34").Value = 2723.25
$ws.Range("L34").Value = 7649.25
$ws.Range("M34").Value = -2639.25
$ws.Range("N34").Value = -7817.25

$ws.Range("H114").Value = 2001
$ws.Range("J114").Value = 2335.3333
$ws.Range("L114").Value = 7005.999899999999
$ws.Range("N114").Value = -13513.9999

$ws.Range("H122").Value = 534
$ws.Range("J122").Value = 800
$ws.Range("L122").Value = 7200
$ws.Range("N122").Value = -12100

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1041.6364
$ws.Range("I107").Value = 1214.1111
$ws.Range("J107").Value = 922.2308
$ws.Range("K107").Value = 1214.1111
$ws.Range("L107").Value = 922.2308
$ws.Range("M107").Value = 705.8888999999999
$ws.Range("N107").Value = -4762.2308

$ws.Range("H113").Value = 5309
$ws.Range("I113").Value = 5495.3335
$ws.Range("J113").Value = 5239.125
$ws.Range("K113").Value = 5495.3335
$ws.Range("L113").Value = 5239.125
$ws.Range("M113").Value = -3325.3335
$ws.Range("N113").Value = -9579.125

$ws.Range("H122").Value = 62016.35
$ws.Range("I122").Value = 2607.6365
$ws.Range("J122").Value = 170932.33
$ws.Range("K122").Value = 7822.9095
$ws.Range("L122").Value = 512796.99
$ws.Range("M122").Value = -5372.9095
$ws.Range("N122").Value = -517696.99

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2621.4
$ws.Range("I7").Value = 2376.75
$ws.Range("J7").Value = 3600
$ws.Range("K7").Value = 2376.75
$ws.Range("L7").Value = 3600
$ws.Range("M7").Value = -2264.75
$ws.Range("N7").Value = -3824

$ws.Range("H16").Value = 1419.1818
$ws.Range("I16").Value = 1419.1818
$ws.Range("K16").Value = 1419.1818
$ws.Range("M16").Value = -1249.1818

$ws.Range("H17").Value = 766.2
$ws.Range("I17").Value = 766.2
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 766.2
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -596.2
$ws.Range("N17").ClearContents()

$ws.Range("H46").Value = 3123.0952
$ws.Range("I46").Value = 2199.8
$ws.Range("J46").Value = 3411.625
$ws.Range("K46").Value = 2199.8
$ws.Range("L46").Value = 3411.625
$ws.Range("M46").Value = -2011.8
$ws.Range("N46").Value = -3787.625

$ws.Range("H93").Value = 1325.091
$ws.Range("I93").Value = 1320.1
$ws.Range("K93").Value = 1320.1
$ws.Range("M93").Value = -72.09999999999991

$ws.Range("H122").Value = 8947.799999999999
$ws.Range("I122").Value = 9184.75
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 27554.25
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -25104.25
$ws.Range("N122").Value = -28900

$ws.Range("H126").Value = 2621.4
$ws.Range("I126").Value = 2376.75
$ws.Range("J126").Value = 3600
$ws.Range("K126").Value = 7130.25
$ws.Range("L126").Value = 10800
$ws.Range("M126").Value = -4660.25
$ws.Range("N126").Value = -15740

$ws.Range("H132").Value = 6519.8
$ws.Range("J132").Value = 13500
$ws.Range("L132").Value = 40500
$ws.Range("N132").Value = -45560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 50019
$ws.Range("J28").Value = 50019
$ws.Range("L28").Value = 50019
$ws.Range("N28").Value = -50715

$ws.Range("H100").Value = 2148.4443
$ws.Range("I100").Value = 1878.4
$ws.Range("K100").Value = 3756.8
$ws.Range("M100").Value = -3215.8

$ws.Range("H126").Value = 2165.818
$ws.Range("I126").Value = 1477.375
$ws.Range("K126").Value = 4432.125
$ws.Range("M126").Value = -1962.125
